# "New BoP pressure drop, equal to architecture 21."
# - L1 header: "P_out in [bar]" -> "bop pressure drop in [bar]"
# - Data table (rows 2..5, index A=0..3) is replaced/extended to rows 2..10
#   (index A=0..8), all numeric values recomputed for the new BoP model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
$ws.Range("L1").Value = "bop pressure drop in [bar]"

# --- Extend the styled index column (col A) down to row 10 -----------
# Rows 6-10 are brand new; copy the formatting (border/bold/center) from
# an existing indexed row so they match rows 2-5.
$ws.Range("A2").Copy()
$ws.Range("A6:A10").PasteSpecial(-4122)

# --- New data table, rows 2..10, columns A..L (index 0..8) -----------
$data = @(
    @(0, 0.1666666666666667, 0.03078718242887758, 0.1358798175711224, 306.6674606968199, 315.1035309651432, 310.9716196524392, 307.6614797046352, 320.3861989161594, 306.6674606968199, 325.6796778876807, 0.07606633230618498),
    @(1, 0.2, 0.03719881950877927, 0.1628011804912188, 309.4205584734395, 316.3973690341517, 312.9793391096121, 310.2422489515129, 320.7688472191118, 309.4205584734395, 325.2580534501679, 0.1014861590692906),
    @(2, 0.25, 0.04684844280810583, 0.2031515571918942, 312.1704476099254, 317.70676633787, 314.9938973253738, 312.8222345175304, 321.1773747330067, 312.1704476099254, 324.836433382633, 0.1459922342955988),
    @(3, 0.2833333333333333, 0.05329556577176932, 0.2300374342282307, 313.4635556567899, 318.3289404428087, 315.9446451588739, 314.036267093806, 321.3795479015832, 313.4635556567899, 324.6380282756148, 0.1799131506506343),
    @(4, 0.3333333333333333, 0.06298031427688508, 0.2703526857231149, 314.9177038178763, 319.0339059454114, 317.0165900071119, 315.4021559487556, 321.6152681676467, 314.9177038178763, 324.4148216346758, 0.2371707588138338),
    @(5, 0.3666666666666666, 0.06944369957187221, 0.2972233004281278, 315.6665988481321, 319.3992706452299, 317.5698539854616, 316.1058821987543, 321.7403163081778, 315.6665988481321, 324.299834870174, 0.2795938541944102),
    @(6, 0.4166666666666667, 0.07914587375499872, 0.3375211262450013, 316.5650743986903, 319.8397726938424, 318.2347508871371, 316.9504326975884, 321.8937813060626, 316.5650743986903, 324.1618541532335, 0.3496029609345089),
    @(7, 0.45, 0.08561766907149518, 0.3643823309285048, 317.0530562460292, 320.0800385932849, 318.5964010967026, 317.4092512500918, 321.9787589386128, 317.0530562460292, 324.0869024983483, 0.4005256201793522),
    @(8, 0.5, 0.09532973264767332, 0.4046702673523267, 317.6629731430928, 320.3813808876369, 319.048958290813, 317.9828440863042, 322.0866340269861, 317.6629731430928, 323.9932120428593, 0.4832860087924997)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
